$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.764.43'
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('D3').Value = '3.072.91'
$ws.Range('E3').Value = '  -1.25%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.51'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '170.49'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '3.071.29'
$ws.Range('E8').Value = '  -1.19%  '
$ws.Range('E9').Value = '  -1.92%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.38'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.15%  '
$ws.Range('E11').Value = '  -1.53%  '
$ws.Range('E12').Value = '  -3.34%  '
$ws.Range('E13').Value = '  -2.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.71'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.89%  '
$ws.Range('E15').Value = '  -1.87%  '
$ws.Range('D16').Value = '3.583.94'
$ws.Range('D17').Value = '66.722.43'
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('E18').Value = '  -2.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.96'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.01%  '
$ws.Range('D20').Value = '3.070.81'
$ws.Range('E20').Value = '  -1.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '489.70'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.71'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.89%  '
$ws.Range('E23').Value = '  -3.65%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.78'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.67'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.24%  '
$ws.Range('E26').Value = '  -3.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.15'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.49%  '
$ws.Range('E28').Value = '  +0.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.81'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.26'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.06%  '
$ws.Range('E31').Value = '  -2.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.53'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.43%  '
$ws.Range('E33').Value = '  -2.61%  '
$ws.Range('D34').Value = '0.0₃0915'
$ws.Range('E34').Value = '  -2.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.951'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.36%  '
$ws.Range('B37').Value = 'Arweave'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '47.26'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.73%  '
$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.58'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.58%  '
$ws.Range('E39').Value = '  +0.30%  '
$ws.Range('E40').Value = '  -4.89%  '
$ws.Range('E41').Value = '  -3.39%  '
$ws.Range('E42').Value = '  -4.56%  '
$ws.Range('D43').Value = '2.756.50'
$ws.Range('E43').Value = '  -2.89%  '
$ws.Range('E44').Value = '  -3.32%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.51'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.66%  '
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '135.36'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '366.46'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.74%  '
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '24.63'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.88%  '
$ws.Range('E50').Value = '  -1.66%  '
$ws.Range('E51').Value = '  -2.00%  '
